$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 22; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value2 = "fullRNASeq"
    }
}
